# Case_3_249/res_bus/vm_pu.xlsx update — "case with 380 kV done"
# Updates vm_pu results for rows 2-25 (bus voltage magnitudes, columns B-F, I-N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=1.02; C=1.036897754395654; D=1.045054397419338; E=1.045713539688314; F=1.056499326006375; I=1.041475486622723; J=1.042003567069894; K=1.047823715058676; L=1.048481006881135; M=1.059236869085619; N=1.01795397991912 }
    3 = @{ B=1.02; C=1.037672824603702; D=1.04565916556913; E=1.04640096635731; F=1.057263585920294; I=1.041657326117372; J=1.042423670028049; K=1.048240384574989; L=1.048980252678307; M=1.059814905820955; N=1.018093652924174 }
    4 = @{ B=1.02; C=1.038175002752569; D=1.046051026858271; E=1.046846756100785; F=1.057759191844736; I=1.041773945629559; J=1.042695465704114; K=1.048509836479969; L=1.04930359643157; M=1.060189347183397; N=1.018183999184924 }
    5 = @{ B=1.02; C=1.038386273537069; D=1.046215891659052; E=1.047034398366809; F=1.057967800603939; I=1.041822721926379; J=1.042809717892125; K=1.048623074036711; L=1.049439599991159; M=1.060346859389729; N=1.018221972723816 }
    6 = @{ B=1.02; C=1.038421755869856; D=1.046243580534343; E=1.047065917888979; F=1.058002841882039; I=1.041830896963012; J=1.042828900669069; K=1.048642084731584; L=1.049462439637148; M=1.060373312021333; N=1.018228348163146 }
    7 = @{ B=1.02; C=1.038177825156153; D=1.04605322929596; E=1.046849262475896; F=1.057761978281944; I=1.041774598366421; J=1.042696992391263; K=1.048511349724832; L=1.049305413445328; M=1.060191451490105; N=1.018184506621369 }
    8 = @{ B=1.02; C=1.037159556031509; D=1.045258669510857; E=1.04594565486715; F=1.056757386555576; I=1.041537155510082; J=1.042145550089142; K=1.047964562966057; L=1.048649666614981; M=1.059432132721566; N=1.018001189298845 }
    9 = @{ B=1.02; C=1.035370344639894; D=1.043862744727283; E=1.044360969259138; F=1.054995526221938; I=1.041110803985039; J=1.041173599969963; K=1.046999884765631; L=1.047496511556336; M=1.058097358330876; N=1.017677940825732 }
    10 = @{ B=1.02; C=1.034181087483996; D=1.042935063782629; E=1.043309735461924; F=1.053826704333329; I=1.040821280281411; J=1.040525554871506; K=1.046356065815736; L=1.046729418910777; M=1.057209796378208; N=1.017462322930539 }
    11 = @{ B=1.02; C=1.033666992548468; D=1.042534088933514; E=1.04285580451462; F=1.053321983721835; I=1.040694670099559; J=1.040244942938943; K=1.0460771373598; L=1.046397677156698; M=1.05682603767871; N=1.017368935926766 }
    12 = @{ B=1.02; C=1.033476166172004; D=1.042385258793457; E=1.042687385806957; F=1.053134718673112; I=1.040647455370388; J=1.040140711955118; K=1.045973509767789; L=1.046274517383727; M=1.056683578823103; N=1.017334244871309 }
    13 = @{ B=1.02; C=1.033517093100875; D=1.042417178363692; E=1.042723503493563; F=1.053174878092393; I=1.04065759149564; J=1.040163069797071; K=1.045995739171089; L=1.046300932656128; M=1.056714132800164; N=1.017341686346644 }
    14 = @{ B=1.02; C=1.033651216089256; D=1.042521784336647; E=1.04284187905377; F=1.05330650002393; I=1.040690771105871; J=1.040236327149243; K=1.046068571893489; L=1.046387495418828; M=1.056814260215839; N=1.017366068410902 }
    15 = @{ B=1.02; C=1.033733871156348; D=1.042586250135644; E=1.0429148395767; F=1.053387624622562; I=1.040711189529698; J=1.040281463567285; K=1.046113443780032; L=1.04644083810778; M=1.056875963500679; N=1.017381090627443 }
    16 = @{ B=1.02; C=1.03421522458633; D=1.04296169044534; E=1.043339888096694; F=1.053860230394911; I=1.040829656838907; J=1.040544178194058; K=1.046374574332263; L=1.046751444388435; M=1.057235277167483; N=1.017468520275521 }
    17 = @{ B=1.02; C=1.034517396855763; D=1.043197387662853; E=1.043606848701341; F=1.054157056454936; I=1.040903635699205; J=1.040708972023731; K=1.046538335355021; L=1.046946391601447; M=1.057460816773505; N=1.017523356767982 }
    18 = @{ B=1.02; C=1.034693731794664; D=1.043334934917891; E=1.043762683756591; F=1.054330323744446; I=1.040946666179263; J=1.040805092984174; K=1.046633839758764; L=1.047060140889698; M=1.057592424260522; N=1.01755533971076 }
    19 = @{ B=1.02; C=1.03475387146432; D=1.043381846628463; E=1.043815840018873; F=1.054389426015946; I=1.040961318056089; J=1.040837867621505; K=1.046666401765733; L=1.047098933171898; M=1.057637308123986; N=1.017566244670055 }
    20 = @{ B=1.02; C=1.034484968035107; D=1.043172092424609; E=1.043578193777378; F=1.054125195985202; I=1.040895710887504; J=1.040691291246209; K=1.046520766849682; L=1.046925471474069; M=1.057436612904169; N=1.017517473561028 }
    21 = @{ B=1.02; C=1.033611716569688; D=1.04249097743829; E=1.04280701509953; F=1.053267734812159; I=1.040681005666835; J=1.040214754652073; K=1.046047125043839; L=1.046362003066165; M=1.056784772792763; N=1.017358888578531 }
    22 = @{ B=1.02; C=1.033063429775654; D=1.04206336952489; E=1.042323253218483; F=1.052729834875187; I=1.040544935875019; J=1.039915142400295; K=1.045749206040332; L=1.046008098536365; M=1.056375434757457; N=1.017259163000907 }
    23 = @{ B=1.02; C=1.033354014178539; D=1.042289991596654; E=1.042579598739036; F=1.053014869412638; I=1.040617170746144; J=1.04007397156121; K=1.045907149572879; L=1.046195674370668; M=1.056592384608773; N=1.017312030853398 }
    24 = @{ B=1.02; C=1.034499620967185; D=1.043183522040854; E=1.043591141329104; F=1.054139591947345; I=1.040899292140113; J=1.040699280430218; K=1.046528705346418; L=1.046934924256138; M=1.057447549424119; N=1.017520131936045 }
    25 = @{ B=1.02; C=1.035832281342802; D=1.044223115934068; E=1.044769736797002; F=1.055450005217074; I=1.041221962291619; J=1.041424892174337; K=1.047249406940274; L=1.04779434102227; M=1.058442034170332; N=1.01776153138091 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

Write-Output "vm_pu.xlsx: case with 380 kV done - updated rows 2-25"